# BIS-1002: removed "Internal Assignment" column from export.
#
# The "Internal Assignment" column (column O) is dropped from the
# exported data-set-type attribute sheet: its header in row 4 and the
# "FALSE" values in rows 5-7 are cleared (cell formatting/style is kept,
# only the content goes away), and the now-unused "Internal Assignment"
# shared string is pruned automatically by the writer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O4:O7").ClearContents()
